$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '78.683.59'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +2.85%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.171.87'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +4.43%  '
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '205.13'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +2.16%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '629.46'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.33%  '
$ws.Range("E7").Value = '  +0.02%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.225'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +10.76%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.579'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +4.94%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '3.169.43'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +4.40%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.576'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +32.54%  '
$ws.Range("E12").Value = '  +2.41%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.40'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +6.45%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '3.754.68'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +4.38%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.0000221'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +16.94%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '31.38'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +6.91%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '78.603.04'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +2.87%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.166.95'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +4.09%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '14.41'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +6.77%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '9.31'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +2.31%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '427.20'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +13.81%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '2.81'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +23.73%  '
$ws.Range("E23").Value = '  +12.78%  '
$ws.Range("E24").Value = '  +5.51%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '3.336.79'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +4.07%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '4.71'
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '75.58'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +3.46%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.84'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +10.28%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.999'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.00%  '
$ws.Range("E30").Value = '  +3.60%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.00'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.32%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '8.82'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +5.88%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.47'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +4.67%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '506.91'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -1.17%  '
$ws.Range("E35").Value = '  +0.36%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.127'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +21.47%  '
$ws.Range("E37").Value = '  +9.25%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.133'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +19.40%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.999'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.03%  '
$ws.Range("B40").Value = 'Monero'
$ws.Range("C40").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '163.96'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.19%  '
$ws.Range("B41").Value = 'PolygonEcosystemToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.395'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +3.41%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '19.94'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.30%  '
$ws.Range("E43").Value = '  -0.16%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '191.19'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.48%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '5.36'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +6.02%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.804'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +14.81%  '
$ws.Range("E47").Value = '  +7.00%  '
$ws.Range("E48").Value = '  +2.86%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '42.56'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.43%  '
$ws.Range("B50").Value = 'dogwifhat'
$ws.Range("C50").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.47'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +4.90%  '
$ws.Range("B51").Value = 'ARBITRUM'
$ws.Range("C51").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.620'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +2.00%  '
